$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("L6")

$columns = @("C", "D", "E")

foreach ($col in $columns) {
    for ($row = 2; $row -le 11; $row++) {
        $cell = $ws.Range("$col$row")
        $text = $cell.Value2

        $parts = $text.Split(",")
        $name = $parts[0]
        $numbersText = $parts[1]

        $numbers = $numbersText.Trim().Split(" ")
        $sum = 0
        foreach ($n in $numbers) {
            if ($n -ne "") {
                $sum += [int]$n
            }
        }

        $newText = "$name,$numbersText,($sum)"
        $cell.Value = $newText
    }
}
